$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.609.28'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '1.632.21'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.10'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.43%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0861'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.52%  '
$ws.Range("D12").Value = '1.857.96'
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").Value = '1.630.40'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.559'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("D17").Value = '27.543.39'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.55%  '
$ws.Range("D19").Value = '0.0₃0719'
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.55'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.29%  '
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.17%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("B29").Value = 'BinanceUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("E31").Value = '  -0.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("D33").Value = '1.467.96'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.68%  '
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.878'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.55%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0167'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.558'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.922'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.39%  '
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '67.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("B45").Value = 'mCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("D47").Value = '1.767.81'
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.74'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0995'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("E51").Value = '  -1.76%  '
